$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'42.834.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  +0.46%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'2.530.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  +0.09%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Formula = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Formula = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'317.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  +1.21%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'96.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  +1.40%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  -1.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Formula = "'  +0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Formula = "'  -0.41%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'35.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  -1.05%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.0820"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  +0.69%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'7.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  -1.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Formula = "'  -3.57%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'2.921.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +0.24%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Formula = "'Chainlink"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Formula = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Formula = "'15.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  -3.84%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Formula = "'WrappedEther"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Formula = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Formula = "'2.447.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  -3.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'0.851"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  -1.22%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'42.913.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Formula = "'  +2.74%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'12.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  -3.45%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'0.0₃0967"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -0.42%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'69.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  -2.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'252.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -0.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Formula = "'2.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Formula = "'  +0.25%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Formula = "'  +1.26%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'26.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  -4.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Formula = "'  +0.80%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Formula = "'  +1.87%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'40.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Formula = "'  +3.14%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Formula = "'10.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  +4.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Formula = "'5.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Formula = "'  -0.35%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'157.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  +1.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'2.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  +0.22%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Formula = "'Celestia"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Formula = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Formula = "'19.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  -4.96%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Formula = "'WEMIXToken"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Formula = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Formula = "'2.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Formula = "'  +3.59%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Formula = "'  -1.41%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'0.0791"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  -0.09%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Formula = "'  +0.24%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Formula = "'  +9.82%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'0.118"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  -0.94%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'22.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  -11.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Formula = "'VeChain"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Formula = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Formula = "'0.0305"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Formula = "'  +0.86%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Formula = "'RenderToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Formula = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Formula = "'3.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  -0.93%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'  +0.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -3.34%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Formula = "'2.010.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Formula = "'  -1.82%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'9.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  +3.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'84.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -1.71%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'106.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  +4.80%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Formula = "'74.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Formula = "'  +0.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'2.777.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  +0.07%  "
$ws.Range("E51").Style = "Normal"
